# Generate Report for Handoff
#
# This refreshes the "Latest Handoff Date" / "Latest Handback DateTime"
# timestamp that was shared by several rows on the Overview, zh-cn and
# de-de sheets. The stale timestamp text is replaced with a freshly
# generated one (per-sheet, since each locale sheet recorded its own
# handback completion time).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$newOverviewTimestamp = "2016-03-21 18:28:59"
$ws.Range("D7").Value  = $newOverviewTimestamp
$ws.Range("D10").Value = $newOverviewTimestamp
$ws.Range("D11").Value = $newOverviewTimestamp
$ws.Range("D12").Value = $newOverviewTimestamp
$ws.Range("D13").Value = $newOverviewTimestamp
$ws.Range("D14").Value = $newOverviewTimestamp
$ws.Range("D15").Value = $newOverviewTimestamp
$ws.Range("D16").Value = $newOverviewTimestamp

# ---- zh-cn sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$newZhCnTimestamp = "2016-03-21 18:28:55"
$ws.Range("E7").Value  = $newZhCnTimestamp
$ws.Range("E10").Value = $newZhCnTimestamp
$ws.Range("E11").Value = $newZhCnTimestamp
$ws.Range("E12").Value = $newZhCnTimestamp
$ws.Range("E13").Value = $newZhCnTimestamp
$ws.Range("E14").Value = $newZhCnTimestamp
$ws.Range("E15").Value = $newZhCnTimestamp
$ws.Range("E16").Value = $newZhCnTimestamp

# ---- de-de sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$newDeDeTimestamp = "2016-03-21 18:28:59"
$ws.Range("E7").Value  = $newDeDeTimestamp
$ws.Range("E10").Value = $newDeDeTimestamp
$ws.Range("E11").Value = $newDeDeTimestamp
$ws.Range("E12").Value = $newDeDeTimestamp
$ws.Range("E13").Value = $newDeDeTimestamp
$ws.Range("E14").Value = $newDeDeTimestamp
$ws.Range("E15").Value = $newDeDeTimestamp
$ws.Range("E16").Value = $newDeDeTimestamp
